$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 41, shifting the existing rows 41..128 down to 42..129
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new data record
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 45014
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112031
$ws.Range("G41").Value = "Poroto verde"
$ws.Range("H41").Value = "Magnum"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 60
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 27000
$ws.Range("M41").Value = 26000
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Provincia de Diguillín"
$ws.Range("P41").Value = 1040
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"
